$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C2").Value = 18.47100257873535
$ws.Range("C3").Value = 17.07315444946289
$ws.Range("C4").Value = 16.80684089660645
$ws.Range("C5").Value = 16.8759822845459
$ws.Range("C6").Value = 17.03095436096191
